$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '237.56'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '21.79'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.390'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05623'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.349'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.7969'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.031'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1389'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07324'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03115'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.02969'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09221'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001667'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.251'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.04770'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0005717'
$ws.Range('E18').Value = '17OneONE'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006230'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.005068'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0001501'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.904'
$ws.Range('E24').Value = '23LEOLEOBestin24h'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04073'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006955'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1039'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003002'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008923'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005442'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6760'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.03722'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002102'
